# Applies numeric corrections to the profit-calculation sheets
# (currentAveragePrice / LevePrice / LeveProfit columns H:N) as produced
# by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1990.2858
$ws.Range("I137").Value = 1528
$ws.Range("K137").Value = 4584
$ws.Range("M137").Value = -2034
$ws.Range("H138").Value = 4901.0415
$ws.Range("I138").Value = 906.4194
$ws.Range("J138").Value = 12185.353
$ws.Range("K138").Value = 2719.2582
$ws.Range("L138").Value = 36556.05899999999
$ws.Range("M138").Value = 2420.7418
$ws.Range("N138").Value = -46836.05899999999

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5980.9165
$ws.Range("I61").Value = 6437.0454
$ws.Range("J61").Value = 963.5
$ws.Range("K61").Value = 6437.0454
$ws.Range("L61").Value = 963.5
$ws.Range("M61").Value = -6225.0454
$ws.Range("N61").Value = -1387.5
$ws.Range("H97").Value = 835.61536
$ws.Range("I97").Value = 657.8570999999999
$ws.Range("J97").Value = 1582.2
$ws.Range("K97").Value = 657.8570999999999
$ws.Range("L97").Value = 1582.2
$ws.Range("M97").Value = -161.8570999999999
$ws.Range("N97").Value = -2574.2
$ws.Range("H122").Value = 1351085.8
$ws.Range("I122").Value = 1351085.8
$ws.Range("K122").Value = 4053257.4
$ws.Range("M122").Value = -4050807.4
$ws.Range("H132").Value = 3548.8333
$ws.Range("I132").Value = 1823.4762
$ws.Range("K132").Value = 5470.4286
$ws.Range("M132").Value = -2940.4286
$ws.Range("H136").Value = 5980.9165
$ws.Range("I136").Value = 6437.0454
$ws.Range("J136").Value = 963.5
$ws.Range("K136").Value = 19311.1362
$ws.Range("L136").Value = 2890.5
$ws.Range("M136").Value = -16761.1362
$ws.Range("N136").Value = -7990.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()  # was -36116
$ws.Range("H96").Value = 40000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 40000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 40000
$ws.Range("M96").ClearContents()  # was -7595.6
$ws.Range("N96").Value = -45492
$ws.Range("H97").Value = 22500
$ws.Range("I97").Value = 4999
$ws.Range("J97").Value = 40001
$ws.Range("K97").Value = 4999
$ws.Range("L97").Value = 40001
$ws.Range("M97").Value = -4008
$ws.Range("N97").Value = -41983
$ws.Range("H99").Value = 100001390
$ws.Range("I99").Value = 111112320
$ws.Range("J99").Value = 3011
$ws.Range("K99").Value = 111112320
$ws.Range("L99").Value = 3011
$ws.Range("M99").Value = -111110822
$ws.Range("N99").Value = -6007
$ws.Range("H100").Value = 61243
$ws.Range("J100").Value = 61243
$ws.Range("L100").Value = 61243
$ws.Range("N100").Value = -63407
$ws.Range("H105").Value = 13570.444
$ws.Range("I105").Value = 19538.545
$ws.Range("K105").Value = 19538.545
$ws.Range("M105").Value = -17791.545

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4303.814
$ws.Range("I31").Value = 966.06665
$ws.Range("J31").Value = 12006.308
$ws.Range("K31").Value = 966.06665
$ws.Range("L31").Value = 12006.308
$ws.Range("M31").Value = -671.06665
$ws.Range("N31").Value = -12596.308
$ws.Range("H34").Value = 4303.814
$ws.Range("I34").Value = 966.06665
$ws.Range("J34").Value = 12006.308
$ws.Range("K34").Value = 966.06665
$ws.Range("L34").Value = 12006.308
$ws.Range("M34").Value = -764.06665
$ws.Range("N34").Value = -12410.308
$ws.Range("H58").Value = 1472.9574
$ws.Range("I58").Value = 868.04
$ws.Range("J58").Value = 2160.3635
$ws.Range("K58").Value = 868.04
$ws.Range("L58").Value = 2160.3635
$ws.Range("M58").Value = -665.04
$ws.Range("N58").Value = -2566.3635
$ws.Range("H80").Value = 26309.334
$ws.Range("J80").Value = 26309.334
$ws.Range("L80").Value = 26309.334
$ws.Range("N80").Value = -28555.334
$ws.Range("H83").Value = 26309.334
$ws.Range("J83").Value = 26309.334
$ws.Range("L83").Value = 78928.00199999999
$ws.Range("N83").Value = -90160.00199999999
$ws.Range("H134").Value = 2424.697
$ws.Range("I134").Value = 2548.1853
$ws.Range("J134").Value = 1869
$ws.Range("K134").Value = 7644.5559
$ws.Range("L134").Value = 5607
$ws.Range("M134").Value = -5109.5559
$ws.Range("N134").Value = -10677
$ws.Range("H136").Value = 1472.9574
$ws.Range("I136").Value = 868.04
$ws.Range("J136").Value = 2160.3635
$ws.Range("K136").Value = 2604.12
$ws.Range("L136").Value = 6481.0905
$ws.Range("M136").Value = -54.11999999999989
$ws.Range("N136").Value = -11581.0905

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 300471.44
$ws.Range("J5").Value = 500375
$ws.Range("L5").Value = 1501125
$ws.Range("N5").Value = -1501349
$ws.Range("H98").Value = 33333968
$ws.Range("J98").Value = 33333968
$ws.Range("L98").Value = 100001904
$ws.Range("N98").Value = -100004900
$ws.Range("H107").Value = 100515.25
$ws.Range("I107").Value = 515.7143
$ws.Range("J107").Value = 154361.16
$ws.Range("K107").Value = 1547.1429
$ws.Range("L107").Value = 463083.48
$ws.Range("M107").Value = 372.8571000000002
$ws.Range("N107").Value = -466923.48
$ws.Range("H113").Value = 333838.4
$ws.Range("I113").Value = 489.5
$ws.Range("J113").Value = 833861.75
$ws.Range("K113").Value = 1468.5
$ws.Range("L113").Value = 2501585.25
$ws.Range("M113").Value = 701.5
$ws.Range("N113").Value = -2505925.25
$ws.Range("H135").Value = 300471.44
$ws.Range("J135").Value = 500375
$ws.Range("L135").Value = 4503375
$ws.Range("N135").Value = -4508445

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 38558.25
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 38558.25
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 38558.25
$ws.Range("M15").ClearContents()  # was -19919
$ws.Range("N15").Value = -39134.25
$ws.Range("H81").Value = 38558.25
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 38558.25
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 38558.25
$ws.Range("M81").ClearContents()  # was -19209
$ws.Range("N81").Value = -40554.25
$ws.Range("H84").Value = 38558.25
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 38558.25
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 115674.75
$ws.Range("M84").ClearContents()  # was -55629
$ws.Range("N84").Value = -125658.75
$ws.Range("H102").Value = 947.2222
$ws.Range("I102").Value = 840.625
$ws.Range("J102").Value = 1800
$ws.Range("K102").Value = 840.625
$ws.Range("L102").Value = 1800
$ws.Range("M102").Value = 781.375
$ws.Range("N102").Value = -5044

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 142860510
$ws.Range("I40").Value = 142860510
$ws.Range("K40").Value = 142860510
$ws.Range("M40").Value = -142860374
$ws.Range("H132").Value = 18063070
$ws.Range("I132").Value = 21674636
$ws.Range("J132").Value = 5244.5
$ws.Range("K132").Value = 65023908
$ws.Range("L132").Value = 15733.5
$ws.Range("M132").Value = -65021378
$ws.Range("N132").Value = -20793.5
$ws.Range("H136").Value = 4090.6743
$ws.Range("J136").Value = 2515.75
$ws.Range("L136").Value = 7547.25
$ws.Range("N136").Value = -12647.25

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1658.0454
$ws.Range("I122").Value = 1499.1111
$ws.Range("J122").Value = 1768.0769
$ws.Range("K122").Value = 4497.3333
$ws.Range("L122").Value = 5304.2307
$ws.Range("M122").Value = -2047.3333
$ws.Range("N122").Value = -10204.2307
$ws.Range("H126").Value = 1291.5385
$ws.Range("I126").Value = 955.7143
$ws.Range("K126").Value = 2867.1429
$ws.Range("M126").Value = -397.1428999999998
$ws.Range("H136").Value = 3067.6296
$ws.Range("I136").Value = 3627.2354
$ws.Range("K136").Value = 10881.7062
$ws.Range("M136").Value = -8331.706200000001
